# Updates the crypto price/volume table to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.937.94'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.292.37'
$ws.Range("E3").Value = '  +1.92%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.70'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '75.09'
$ws.Range("E7").Value = '  +7.15%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.647'
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.28'
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.54'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").Value = '2.634.76'
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.07'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.871'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").Value = '2.297.84'
$ws.Range("E17").Value = '  +2.23%  '
$ws.Range("D18").Value = '42.820.47'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  +1.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.25'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.37'
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.99'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +5.49%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.89'
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("E28").Value = '  -2.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.39'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.18'
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0862'
$ws.Range("E31").Value = '  +8.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.26'
$ws.Range("E32").Value = '  -5.38%  '
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.20'
$ws.Range("E34").Value = '  +0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.128'
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.65'
$ws.Range("E36").Value = '  +5.99%  '
$ws.Range("E37").Value = '  +2.04%  '
$ws.Range("E38").Value = '  -4.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.60'
$ws.Range("E39").Value = '  +9.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.96'
$ws.Range("E41").Value = '  +2.32%  '
$ws.Range("E42").Value = '  +4.43%  '
$ws.Range("E43").Value = '  +1.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.21'
$ws.Range("E44").Value = '  -2.41%  '
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '105.13'
$ws.Range("E46").Value = '  +10.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.102'
$ws.Range("E47").Value = '  -1.42%  '
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.17'
$ws.Range("E49").Value = '  -0.73%  '
$ws.Range("E50").Value = '  -1.43%  '
$ws.Range("E51").Value = '  -1.88%  '
